$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.105.01'

# Row 3
$ws.Range("D3").Value = '3.510.88'
$ws.Range("E3").Value = '  +0.49%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = "'594.47"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.87%  '

# Row 6
$ws.Range("D6").Value = "'173.47"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +1.87%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +2.32%  '

# Row 9
$ws.Range("D9").Value = "'0.132"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +5.80%  '

# Row 10
$ws.Range("D10").Value = "'7.28"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -0.22%  '

# Row 11
$ws.Range("D11").Value = "'0.435"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.32%  '

# Row 12
$ws.Range("D12").Value = '4.117.68'
$ws.Range("E12").Value = '  +0.51%  '

# Row 13
$ws.Range("E13").Value = '  +0.15%  '

# Row 14
$ws.Range("D14").Value = "'28.95"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +2.78%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '67.042.79'
$ws.Range("E15").Value = '  +0.70%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +1.45%  '

# Row 17
$ws.Range("D17").Value = '3.487.16'
$ws.Range("E17").Value = '  -0.02%  '

# Row 18
$ws.Range("E18").Value = '  -0.14%  '

# Row 19
$ws.Range("D19").Value = "'14.20"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +2.04%  '

# Row 20
$ws.Range("D20").Value = "'393.75"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +1.42%  '

# Row 21
$ws.Range("D21").Value = "'7.98"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.24%  '

# Row 22
$ws.Range("D22").Value = "'73.12"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.37%  '

# Row 23
$ws.Range("D23").Value = "'0.540"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.72%  '

# Row 24
$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("E25").Value = '  +0.05%  '

# Row 26
$ws.Range("D26").Value = "'10.25"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.61%  '

# Row 27
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("D29").Value = "'6.28"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.37%  '

# Row 30
$ws.Range("D30").Value = "'1.45"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -2.35%  '

# Row 31
$ws.Range("E31").Value = '  +0.88%  '

# Row 32
$ws.Range("D32").Value = "'23.87"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +1.60%  '

# Row 33
$ws.Range("E33").Value = '  -0.64%  '

# Row 34
$ws.Range("E34").Value = '  +1.94%  '

# Row 35
$ws.Range("D35").Value = "'163.28"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +0.53%  '

# Row 36
$ws.Range("D36").Value = "'0.894"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +1.93%  '

# Row 37
$ws.Range("D37").Value = "'1.90"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -0.51%  '

# Row 38
$ws.Range("D38").Value = "'7.08"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +7.38%  '

# Row 39
$ws.Range("D39").Value = "'4.69"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +0.30%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.0748"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +0.05%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = "'27.42"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +3.44%  '

# Row 42
$ws.Range("D42").Value = "'26.41"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.86%  '

# Row 43
$ws.Range("E43").Value = '  +4.67%  '

# Row 44
$ws.Range("D44").Value = '2.815.86'
$ws.Range("E44").Value = '  +0.13%  '

# Row 45
$ws.Range("D45").Value = "'42.99"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -0.67%  '

# Row 46
$ws.Range("D46").Value = "'0.0306"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.21%  '

# Row 47
$ws.Range("D47").Value = "'338.36"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -4.58%  '

# Row 48
$ws.Range("D48").Value = "'1.08"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +0.05%  '

# Row 49
$ws.Range("D49").Value = "'34.19"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +0.97%  '

# Row 50
$ws.Range("D50").Value = "'6.50"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.39%  '

# Row 51
$ws.Range("D51").Value = "'0.848"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -0.50%  '
